# Insert a new weekly price record as row 60 in the "Jengibre" sheet.
# All the existing rows from 60 downward shift down by one row
# (old row 60 becomes row 61, old row 124 becomes row 125, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 60, pushing the rest down.
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new record's data.
$ws.Range("A60").Value = 6
$ws.Range("B60").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C60").Value = "Metropolitana"
$ws.Range("D60").Value = 45040
$ws.Range("E60").Value = 13
$ws.Range("F60").Value = 100114007
$ws.Range("G60").Value = "Jengibre"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 580
$ws.Range("K60").Value = 14000
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = 14448
$ws.Range("N60").Value = "$/caja 13 kilos"
$ws.Range("O60").Value = "Perú"
$ws.Range("P60").Value = 1111
$ws.Range("Q60").Value = 13
$ws.Range("R60").Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D60").NumberFormat = $ws.Range("D61").NumberFormat
